$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: id=2, code=A02, description=test
# Leading apostrophe forces "2" to be stored as text (matching the existing
# row 2 where "1" is also stored as a text shared string) instead of being
# auto-converted to a numeric value.
$ws.Range("A3").Value = "'2"
$ws.Range("B3").Value = "A02"
$ws.Range("C3").Value = "test"
